$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "ExportFile" API row (row 5, columns B and C)
$ws.Range("B5").Value = "ExportFile"
$ws.Range("C5").Value = "Export data ra file xlsx, lưu trên thư mục ExportFiles (temp)"

# Move the active selection to C6, matching the saved selection state
$ws.Range("C6").Select()
